$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date serial for every
# data row (2-439). The sheet was refreshed by one day: 46074 -> 46075.
$ws.Range("C2:C439").Value = 46075
